# Added new Config option: Enable Dining Spot Sign.
# Appends two new rows (toggle09 / tooltip13) to the translations sheet,
# matching the existing id/en/jp/cn table layout, and nudges the sheet's
# view state (row 25 height, active selection) to match the authored diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 30: toggle09 / Enable Dining Spot Sign ---------------------
$ws.Cells.Item(30, 1).Value2 = "toggle09"
$ws.Cells.Item(30, 2).Value2 = "Enable Dining Spot Sign"
$ws.Cells.Item(30, 3).Value2 = "食堂の立て札を有効化"
$ws.Cells.Item(30, 4).Value2 = "启用食堂招牌"

# --- New row 31: tooltip13 / Enable or disable dining spot sign ... -----
$ws.Cells.Item(31, 1).Value2 = "tooltip13"
$ws.Cells.Item(31, 2).Value2 = "Enable or disable dining spot sign effects inside tents."
$ws.Cells.Item(31, 3).Value2 = "テント内で食堂の立て札の効果を有効または無効にします。"
$ws.Cells.Item(31, 4).Value2 = "启用或禁用帐篷内的食堂招牌效果。"

# Columns C/D on the existing data rows use a wrapped "Noto Sans SC" style
# (s="4" / fontId 5). Give the two new rows the same look-and-feel so they
# render consistently with the rest of the sheet.
foreach ($r in 30, 31) {
    $ws.Cells.Item($r, 3).Font.Name = "Noto Sans SC"
    $ws.Cells.Item($r, 3).Font.Size = 10
    $ws.Cells.Item($r, 3).WrapText = $true

    $ws.Cells.Item($r, 4).Font.Name = "Noto Sans SC"
    $ws.Cells.Item($r, 4).Font.Size = 10
    $ws.Cells.Item($r, 4).WrapText = $true

    # Matches the row height ("ht=12.8", customHeight) used by the new rows.
    $ws.Rows.Item($r).RowHeight = 12.8
}

# Row 25 (tooltip10) height was also adjusted as part of this commit.
$ws.Rows.Item(25).RowHeight = 13.4

# Move the sheet's active selection the way it ended up after the edit.
[void]$ws.Range("C42").Select()
